$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1277.6364
$ws.Range("I32").Value = 1100
$ws.Range("J32").Value = 1344.25
$ws.Range("K32").Value = 1100
$ws.Range("L32").Value = 1344.25
$ws.Range("M32").Value = -774
$ws.Range("N32").Value = -1996.25

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 4790507
$ws.Range("I86").Value = 7499.75
$ws.Range("J86").Value = 7523653.5
$ws.Range("K86").Value = 7499.75
$ws.Range("L86").Value = 7523653.5
$ws.Range("M86").Value = -6376.75
$ws.Range("N86").Value = -7525899.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 4790507
$ws.Range("I89").Value = 7499.75
$ws.Range("J89").Value = 7523653.5
$ws.Range("K89").Value = 37498.75
$ws.Range("L89").Value = 37618267.5
$ws.Range("M89").Value = -31882.75
$ws.Range("N89").Value = -37629499.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H101").Value = 856
$ws.Range("I101").Value = 840
$ws.Range("K101").Value = 2520
$ws.Range("M101").Value = -898

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 75585.78999999999
$ws.Range("I111").Value = 86851
$ws.Range("J111").Value = 7994.5
$ws.Range("K111").Value = 260553
$ws.Range("L111").Value = 23983.5
$ws.Range("M111").Value = -257486
$ws.Range("N111").Value = -30117.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 11350.8125
$ws.Range("I116").Value = 4561.778
$ws.Range("J116").Value = 20079.572
$ws.Range("K116").Value = 4561.778
$ws.Range("L116").Value = 20079.572
$ws.Range("M116").Value = -1119.778
$ws.Range("N116").Value = -26963.572

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 1886.8518
$ws.Range("I129").Value = 697.5833
$ws.Range("K129").Value = 2092.7499
$ws.Range("M129").Value = 2907.2501

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H131").Value = 4849.8184
$ws.Range("I131").Value = 3731
$ws.Range("K131").Value = 11193
$ws.Range("M131").Value = -6153

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 4309.1333
$ws.Range("I132").Value = 4472.077
$ws.Range("K132").Value = 13416.231
$ws.Range("M132").Value = -10886.231

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 805742.4399999999
$ws.Range("I137").Value = 627450.5
$ws.Range("K137").Value = 1882351.5
$ws.Range("M137").Value = -1879801.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 5905.552
$ws.Range("J138").Value = 6674
$ws.Range("L138").Value = 20022
$ws.Range("N138").Value = -30302

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5357.537
$ws.Range("I32").Value = 4703.906
$ws.Range("K32").Value = 4703.906
$ws.Range("M32").Value = -4416.906

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2573.2173
$ws.Range("I61").Value = 1936.7
$ws.Range("K61").Value = 1936.7
$ws.Range("M61").Value = -1724.7

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2263.9048
$ws.Range("I74").Value = 2177.15
$ws.Range("K74").Value = 2177.15
$ws.Range("M74").Value = -1303.15

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 2263.9048
$ws.Range("I77").Value = 2177.15
$ws.Range("K77").Value = 10885.75
$ws.Range("M77").Value = -6517.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2573.2173
$ws.Range("I136").Value = 1936.7
$ws.Range("K136").Value = 5810.1
$ws.Range("M136").Value = -3260.1

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H21").Value = 39999
$ws.Range("J21").Value = 39999
$ws.Range("L21").Value = 39999
$ws.Range("N21").Value = -40471

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1194.3
$ws.Range("I105").Value = 1143
$ws.Range("J105").Value = 1399.5
$ws.Range("K105").Value = 1143
$ws.Range("L105").Value = 1399.5
$ws.Range("M105").Value = 604
$ws.Range("N105").Value = -4893.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 62729.723
$ws.Range("I31").Value = 2917.4
$ws.Range("K31").Value = 2917.4
$ws.Range("M31").Value = -2622.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 62729.723
$ws.Range("I34").Value = 2917.4
$ws.Range("K34").Value = 2917.4
$ws.Range("M34").Value = -2715.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H92").Value = 100000
$ws.Range("J92").Value = 100000
$ws.Range("L92").Value = 100000
$ws.Range("N92").Value = -104992

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 151088.67
$ws.Range("I99").Value = 3538.6316
$ws.Range("J99").Value = 337985.4
$ws.Range("K99").Value = 3538.6316
$ws.Range("L99").Value = 337985.4
$ws.Range("M99").Value = -2040.6316
$ws.Range("N99").Value = -340981.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 2705.5806
$ws.Range("I122").Value = 1519.2
$ws.Range("J122").Value = 4862.636
$ws.Range("K122").Value = 4557.6
$ws.Range("L122").Value = 14587.908
$ws.Range("M122").Value = -2107.6
$ws.Range("N122").Value = -19487.908

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 151088.67
$ws.Range("I126").Value = 3538.6316
$ws.Range("J126").Value = 337985.4
$ws.Range("K126").Value = 10615.8948
$ws.Range("L126").Value = 1013956.2
$ws.Range("M126").Value = -8145.8948
$ws.Range("N126").Value = -1018896.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 4154.9688
$ws.Range("I132").Value = 3662.3333
$ws.Range("J132").Value = 5632.875
$ws.Range("K132").Value = 10986.9999
$ws.Range("L132").Value = 16898.625
$ws.Range("M132").Value = -8456.999899999999
$ws.Range("N132").Value = -21958.625

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1115595.5
$ws.Range("I134").Value = 773133.7
$ws.Range("K134").Value = 2319401.1
$ws.Range("M134").Value = -2316866.1

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 4754749
$ws.Range("J33").Value = 200151.2
$ws.Range("L33").Value = 1200907.2
$ws.Range("N33").Value = -1201473.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 3334363.2
$ws.Range("I68").Value = 10000002
$ws.Range("J68").Value = 2001235.6
$ws.Range("K68").Value = 30000006
$ws.Range("L68").Value = 6003706.800000001
$ws.Range("M68").Value = -29999195
$ws.Range("N68").Value = -6005328.800000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 3334363.2
$ws.Range("I71").Value = 10000002
$ws.Range("J71").Value = 2001235.6
$ws.Range("K71").Value = 90000018
$ws.Range("L71").Value = 18011120.4
$ws.Range("M71").Value = -89995962
$ws.Range("N71").Value = -18019232.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 3209.8333
$ws.Range("I98").Value = 3407.5
$ws.Range("J98").Value = 3111
$ws.Range("K98").Value = 10222.5
$ws.Range("L98").Value = 9333
$ws.Range("M98").Value = -8724.5
$ws.Range("N98").Value = -12329

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 48623.047
$ws.Range("I107").Value = 978.0909
$ws.Range("J107").Value = 96268
$ws.Range("K107").Value = 2934.2727
$ws.Range("L107").Value = 288804
$ws.Range("M107").Value = -1014.2727
$ws.Range("N107").Value = -292644

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 111113400
$ws.Range("J114").Value = 2576.375
$ws.Range("L114").Value = 7729.125
$ws.Range("N114").Value = -14237.125

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 20196670
$ws.Range("I121").Value = 541.55554
$ws.Range("J121").Value = 42917316
$ws.Range("K121").Value = 1624.66662
$ws.Range("L121").Value = 128751948
$ws.Range("M121").Value = -314.66662
$ws.Range("N121").Value = -128754568

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 1594.8
$ws.Range("J129").Value = 2305.5
$ws.Range("L129").Value = 6916.5
$ws.Range("N129").Value = -16916.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 1999.6957
$ws.Range("I137").Value = 1146.5294
$ws.Range("J137").Value = 4417
$ws.Range("K137").Value = 3439.5882
$ws.Range("L137").Value = 13251
$ws.Range("M137").Value = 1660.4118
$ws.Range("N137").Value = -23451

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H139").Value = 6424.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 1821527.6
$ws.Range("J80").Value = 2001961.4
$ws.Range("L80").Value = 2001961.4
$ws.Range("N80").Value = -2003957.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 1821527.6
$ws.Range("J83").Value = 2001961.4
$ws.Range("L83").Value = 10009807
$ws.Range("N83").Value = -10019791

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2351.3333
$ws.Range("I102").Value = 1558.04
$ws.Range("K102").Value = 1558.04
$ws.Range("M102").Value = 63.96000000000004

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 564247.7
$ws.Range("I113").Value = 911791.6
$ws.Range("K113").Value = 911791.6
$ws.Range("M113").Value = -909621.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1868.3334
$ws.Range("I122").Value = 1239.5333
$ws.Range("J122").Value = 3125.9333
$ws.Range("K122").Value = 3718.5999
$ws.Range("L122").Value = 9377.7999
$ws.Range("M122").Value = -1268.5999
$ws.Range("N122").Value = -14277.7999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 359392.6
$ws.Range("I132").Value = 421563.4
$ws.Range("J132").Value = 146235.42
$ws.Range("K132").Value = 1264690.2
$ws.Range("L132").Value = 438706.26
$ws.Range("M132").Value = -1262160.2
$ws.Range("N132").Value = -443766.26

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 997.75
$ws.Range("I22").Value = 995
$ws.Range("J22").Value = 998.1429000000001
$ws.Range("K22").Value = 995
$ws.Range("L22").Value = 998.1429000000001
$ws.Range("M22").Value = -700
$ws.Range("N22").Value = -1588.1429

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 997.75
$ws.Range("I27").Value = 995
$ws.Range("J27").Value = 998.1429000000001
$ws.Range("K27").Value = 995
$ws.Range("L27").Value = 998.1429000000001
$ws.Range("M27").Value = -888
$ws.Range("N27").Value = -1212.1429

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 29360.75
$ws.Range("I2").Value = 29299.7
$ws.Range("K2").Value = 29299.7
$ws.Range("M2").Value = -29187.7

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 25644826
$ws.Range("I122").Value = 32261388
$ws.Range("K122").Value = 96784164
$ws.Range("M122").Value = -96781714

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 31046.277
$ws.Range("I132").Value = 2760.1155
$ws.Range("J132").Value = 104590.3
$ws.Range("K132").Value = 8280.3465
$ws.Range("L132").Value = 313770.9
$ws.Range("M132").Value = -5750.3465
$ws.Range("N132").Value = -318830.9

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 11527394
$ws.Range("I136").Value = 19098290
$ws.Range("J136").Value = 171049.67
$ws.Range("K136").Value = 57294870
$ws.Range("L136").Value = 513149.01
$ws.Range("M136").Value = -57292320
